$d = $word.ActiveDocument

# Locate the "User story 2:" heading paragraph, then the body paragraph that
# immediately follows it (this is the paragraph that currently holds the
# trailing _GoBack bookmark).
$story2BodyIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^User story 2:") {
        $story2BodyIndex = $i + 1
        break
    }
}

$story2Body = $d.Paragraphs.Item($story2BodyIndex)

# Insert a new (still empty) paragraph right after the "User story 2" body
# paragraph; it inherits the (non-bold) run formatting from the end of that
# paragraph, so fill it in with the "User story 3" body text first.
$story2Body.Range.InsertParagraphAfter()
$story3Body = $d.Paragraphs.Item($story2BodyIndex + 1)
$story3Body.Range.Text = "For the user, I would like to add a save button to save the layout and customization of a city. Priority is High and estimate time is 2 days. For testing, press the save button and save the file of the city with name."

# Insert the "User story 3:" heading paragraph between the "User story 2"
# body paragraph and the body paragraph we just created, then bold it.
$story2Body = $d.Paragraphs.Item($story2BodyIndex)
$story2Body.Range.InsertParagraphAfter()
$story3Heading = $d.Paragraphs.Item($story2BodyIndex + 1)
$story3Heading.Range.Text = "User story 3:"
$story3Heading.Range.Bold = 1

# Move the _GoBack bookmark from the "User story 2" body paragraph to the
# (now last) trailing empty paragraph, matching Word's behaviour of keeping
# that bookmark at the last edit position in the document.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$lastParaRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$collapsedEnd = $lastParaRange.Duplicate
$collapsedEnd.Collapse(1)
$d.Bookmarks.Add("_GoBack", $collapsedEnd)
